$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 3 ("Summary of Changes" -> "Details of Changes", effort 0 -> 10.5) ---
$ws.Range("C3").Value = "Details of Changes"
$ws.Range("E3").Value = 10.5

# --- Update row 4 (add Items text, effort 0 -> 1) ---
$ws.Range("D4").Value = "Static Data Configurations Impact"
$ws.Range("E4").Value = 1

# --- Remove the merged block C4:C6 before deleting the now-unused rows 5-7 ---
$ws.Range("C4:C6").UnMerge()

# Rows 5 ("Summary" helper row 3), 6 (helper row 4) and 7 ("Others") are no
# longer part of the breakdown table; delete them so the SUM row shifts up
# from row 8 to row 5.
$ws.Range("A5:E7").EntireRow.Delete()

# --- Fix up the total formula so it sums the Effort column directly above it ---
$ws.Range("E5").Formula = "=SUM(E2:E4)"

# --- Column widths: Category column widens, Items column narrows ---
# (Excel's ColumnWidth property is offset from the stored XML "width" by the
# workbook's default ~5px padding, i.e. ~0.8333 characters at this font, so
# back that off here to land on the exact target width attribute values.)
$ws.Columns("C").ColumnWidth = 28 - 0.8333333333333333
$ws.Columns("D").ColumnWidth = 70 - 0.8333333333333333
